# Edit: "run once when there is no testData existed for the test case"
#
# 1. TestCases sheet: re-order rows so that the test case that has NO
#    matching TestData ("verifyDownload") is moved to the bottom of the
#    list (it now "runs once" instead of being interleaved with the ones
#    that do have TestData rows).
# 2. TestData sheet: the TestData row order for loginMY/loginEN is
#    swapped, and the placeholder rows for verifyDownload / deposit / close
#    (test cases that never had real TestData beyond the TestCaseID) are
#    removed entirely.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) TestCases sheet - reorder so "verifyDownload" (no RunMode/Results)
#    moves from row 3 down to the last row (row 5).
# ---------------------------------------------------------------------
$tc = $wb.Worksheets.Item("TestCases")

# capture current row 3 (verifyDownload) before overwriting it
$vdA = $tc.Cells.Item(3, 1).Value()
$vdB = $tc.Cells.Item(3, 2).Value()
$vdC = $tc.Cells.Item(3, 3).Value()

# shift row4 (deposit) up into row3
$tc.Cells.Item(3, 1).Value = $tc.Cells.Item(4, 1).Value()
$tc.Cells.Item(3, 2).Value = $tc.Cells.Item(4, 2).Value()
$tc.Cells.Item(3, 3).Value = $tc.Cells.Item(4, 3).Value()
$tc.Cells.Item(3, 4).Value = $tc.Cells.Item(4, 4).Value()

# shift row5 (close) up into row4
$tc.Cells.Item(4, 1).Value = $tc.Cells.Item(5, 1).Value()
$tc.Cells.Item(4, 2).Value = $tc.Cells.Item(5, 2).Value()
$tc.Cells.Item(4, 3).Value = $tc.Cells.Item(5, 3).Value()
$tc.Cells.Item(4, 4).Value = $tc.Cells.Item(5, 4).Value()

# put the saved verifyDownload row into row5, with no Results value
$tc.Cells.Item(5, 1).Value = $vdA
$tc.Cells.Item(5, 2).Value = $vdB
$tc.Cells.Item(5, 3).Value = $vdC
$tc.Cells.Item(5, 4).ClearContents()

# ---------------------------------------------------------------------
# 2) TestData sheet - swap loginEN/loginMY rows, drop the rows that only
#    had a TestCaseID and no actual test data (verifyDownload, deposit,
#    close) since those test cases run once without TestData.
# ---------------------------------------------------------------------
$td = $wb.Worksheets.Item("TestData")

$row2 = $td.Cells.Item(2, 1).Value()
$row3 = $td.Cells.Item(3, 1).Value()
$td.Cells.Item(2, 1).Value = $row3
$td.Cells.Item(3, 1).Value = $row2

$td.Range("A4:D6").ClearContents()
$td.Range("4:6").Delete()

# ---------------------------------------------------------------------
# 3) Restore view/selection state to match the post-edit workbook.
# ---------------------------------------------------------------------
$ts = $wb.Worksheets.Item("TestSteps")
$ts.Range("G19").Select()

$td.Range("A4").Select()

$tc.Activate()
$tc.Range("B11").Select()

$wb.Save()
